$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = "Appointment Date : 16/11/2023, Time : [ 09:00 AM to 09:04 AM ]"
$ws.Range("AC2").Value = "16/11/2023"
$ws.Range("AS2").Value = "voice_record_16112023"
$ws.Range("AV2").Value = "formshow_16112023"
